$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlCenter = -4108

# --- New column widths (G and I) ---
$ws.Columns.Item(7).ColumnWidth = 15.592447916666666
$ws.Columns.Item(9).ColumnWidth = 17.166666666666668

# --- Row 1 : new "TABELS" header row ---
$ws.Range("B1").Value = "TABELS"
$ws.Range("B1").HorizontalAlignment = $xlCenter

$ws.Range("D1:F1").Value = "NAMA TABEL"
$ws.Range("D1:F1").Merge()
$ws.Range("D1:F1").HorizontalAlignment = $xlCenter

$ws.Range("G1").Value = "TYPE"
$ws.Range("I1").Value = "CO"

# --- Column G / I annotations next to each existing field ---
$ws.Range("G4").Value = "INT25"

$ws.Range("G5").Value = "VARCHAR60"
$ws.Range("I5").Value = "UTF8_general_ci"

$ws.Range("G6").Value = "VARCHAR30"
$ws.Range("I6").Value = "UTF8_general_ci"

$ws.Range("G7").Value = "VARCHAR30"
$ws.Range("I7").Value = "UTF8_general_ci"

$ws.Range("G8").Value = "VARCHAR30"
$ws.Range("I8").Value = "UTF8_general_ci"

$ws.Range("G9").Value = "INT11"

$ws.Range("G10").Value = "INT11"

$ws.Range("G13").Value = "INT25"

$ws.Range("G14").Value = "TEXT"
$ws.Range("I14").Value = "UTF8_general_ci"

$ws.Range("G15").Value = "TEXT"
$ws.Range("I15").Value = "UTF8_general_ci"

$ws.Range("G16").Value = "VARCHAR100"
$ws.Range("I16").Value = "UTF8_general_ci"

$ws.Range("G17").Value = "INT15"

$ws.Range("G21").Value = "INT25"

$ws.Range("G22").Value = "INT11"

$ws.Range("G23").Value = "LONGTEXT"
$ws.Range("I23").Value = "UTF8_general_ci"

# --- View tweaks ---
$ws.Application.ActiveWindow.DisplayGridlines = $true
$ws.Range("I23").Select()

Write-Host "edit applied"
